# Apply crypto price/volume updates (GitHub Actions data refresh, 2023-05-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold text that looks numeric (e.g. "91.60", "26.975.43").
# Force text format first so Excel does not re-parse them as numbers and drop
# significant trailing zeros / thousands-style dots.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") "26.975.43"
$ws.Range("E2").Value = "  +0.16%  "

Set-TextValue $ws.Range("D3") "1.819.45"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.23%  "

Set-TextValue $ws.Range("D5") "310.38"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("E6").Value = "  +0.22%  "

Set-TextValue $ws.Range("D7") "0.4679"
$ws.Range("E7").Value = "  +0.61%  "

Set-TextValue $ws.Range("D8") "0.3667"
$ws.Range("E8").Value = "  -1.05%  "

Set-TextValue $ws.Range("D9") "0.07350"
$ws.Range("E9").Value = "  -0.27%  "

Set-TextValue $ws.Range("D10") "0.8735"
$ws.Range("E10").Value = "  -0.02%  "

Set-TextValue $ws.Range("D11") "20.28"
$ws.Range("E11").Value = "  -1.14%  "

Set-TextValue $ws.Range("D12") "1.836.47"
$ws.Range("E12").Value = "  -0.21%  "

Set-TextValue $ws.Range("D13") "5.418"
$ws.Range("E13").Value = "  +1.22%  "

Set-TextValue $ws.Range("D14") "0.07120"
$ws.Range("E14").Value = "  +0.69%  "

Set-TextValue $ws.Range("D15") "6.514"
$ws.Range("E15").Value = "  -0.20%  "

Set-TextValue $ws.Range("D16") "91.60"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("E17").Value = "  +0.24%  "

Set-TextValue $ws.Range("D18") "0.000008732"
$ws.Range("E18").Value = "  -0.15%  "

Set-TextValue $ws.Range("D20") "14.69"
$ws.Range("E20").Value = "  -0.46%  "

Set-TextValue $ws.Range("D21") "27.005.06"
$ws.Range("E21").Value = "  +0.20%  "

Set-TextValue $ws.Range("D22") "5.295"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("E23").Value = "  +0.01%  "

Set-TextValue $ws.Range("D24") "2.038.09"
$ws.Range("E24").Value = "  -0.58%  "

Set-TextValue $ws.Range("D25") "1.895"
$ws.Range("E25").Value = "  -0.70%  "

Set-TextValue $ws.Range("D26") "151.09"
$ws.Range("E26").Value = "  -0.59%  "

Set-TextValue $ws.Range("D27") "18.40"
$ws.Range("E27").Value = "  -0.04%  "

Set-TextValue $ws.Range("D28") "2.142"
$ws.Range("E28").Value = "  -0.72%  "

Set-TextValue $ws.Range("D29") "5.243"
$ws.Range("E29").Value = "  -1.70%  "

Set-TextValue $ws.Range("D30") "116.91"
$ws.Range("E30").Value = "  +1.00%  "

Set-TextValue $ws.Range("D31") "0.08859"
$ws.Range("E31").Value = "  -0.54%  "

Set-TextValue $ws.Range("D32") "0.7581"
$ws.Range("E32").Value = "  -0.21%  "

Set-TextValue $ws.Range("D33") "1.160"
$ws.Range("E33").Value = "  +0.30%  "

Set-TextValue $ws.Range("D34") "4.502"
$ws.Range("E34").Value = "  +0.58%  "

Set-TextValue $ws.Range("D35") "2.932"
$ws.Range("E35").Value = "  +0.33%  "

$ws.Range("E36").Value = "  +0.23%  "

Set-TextValue $ws.Range("D37") "1.101"
$ws.Range("E37").Value = "  +0.61%  "

Set-TextValue $ws.Range("D38") "0.05309"
$ws.Range("E38").Value = "  +0.79%  "

Set-TextValue $ws.Range("D39") "0.01947"
$ws.Range("E39").Value = "  -0.66%  "

Set-TextValue $ws.Range("D40") "2.975"
$ws.Range("E40").Value = "  +1.29%  "

Set-TextValue $ws.Range("D44") "0.1652"
$ws.Range("E44").Value = "  -0.76%  "

Set-TextValue $ws.Range("D45") "8.456"
$ws.Range("E45").Value = "  -0.16%  "

Set-TextValue $ws.Range("D46") "0.4896"
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("E48").Value = "  +0.24%  "

Set-TextValue $ws.Range("D49") "1.663"
$ws.Range("E49").Value = "  -1.36%  "

Set-TextValue $ws.Range("D50") "103.23"
$ws.Range("E50").Value = "  -0.19%  "

Set-TextValue $ws.Range("D51") "0.06294"
$ws.Range("E51").Value = "  +0.01%  "

# Rows 41-43: coin ranking reshuffled (cyclic rotation) with refreshed price/volume
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D41") "2.368"
$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.5299"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "7.177"
$ws.Range("E43").Value = "  -1.31%  "

